$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, pushing existing rows 34-37 down to 35-38.
$ws.Rows.Item(34).Insert(-4121) # xlShiftDown

# Seed the new row's formatting from row 33 (same TCID/OP/result column layout).
$ws.Range("A33:E33").Copy($ws.Range("A34:E34"))

# Set the new script's values.
$ws.Range("A34").Value = "DRAIAM114"
$ws.Range("B34").Value = "OP114"
$ws.Range("C34").Value = "New User"

# Update the view to mirror the authored state (new row selected/in view).
$ws.Range("C34").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 26
